$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($ws, $addr) {
    # NOTE: deleting a Hyperlinks item while mid-enumeration of the same
    # collection corrupts the underlying iterator (it silently skips the
    # element following a deleted one), so every removal gets its own
    # fresh enumeration that stops as soon as it finds + deletes the match.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------
# Sheet "Overview": drop the row for d26906ea-...-md (row 3) entirely,
# including its hyperlink.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Remove-HyperlinkAt $wsOverview '$A$3'
$wsOverview.Rows.Item(3).Delete()

# ---------------------------------------------------------------------
# Sheet "zh-cn": drop row 3 (the d26906ea... handback) and its
# hyperlinks, then refresh the Correspond Handoff/Handback datetimes for
# the remaining row (row 2).
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Remove-HyperlinkAt $wsZhCn '$A$3'
Remove-HyperlinkAt $wsZhCn '$B$3'
Remove-HyperlinkAt $wsZhCn '$D$3'
Remove-HyperlinkAt $wsZhCn '$F$3'
Remove-HyperlinkAt $wsZhCn '$G$3'
$wsZhCn.Rows.Item(3).Delete()

$wsZhCn.Range("E2").Value = "2016-03-18 10:36:25"
$wsZhCn.Range("H2").Value = "2016-03-18 10:36:43"

# ---------------------------------------------------------------------
# Sheet "de-de": drop row 3 (the d26906ea... handback) and its
# hyperlinks, then refresh the Correspond Handoff/Handback datetimes for
# the remaining row (row 2).
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Remove-HyperlinkAt $wsDeDe '$A$3'
Remove-HyperlinkAt $wsDeDe '$B$3'
Remove-HyperlinkAt $wsDeDe '$D$3'
Remove-HyperlinkAt $wsDeDe '$F$3'
Remove-HyperlinkAt $wsDeDe '$G$3'
$wsDeDe.Rows.Item(3).Delete()

$wsDeDe.Range("E2").Value = "2016-03-18 10:36:27"
$wsDeDe.Range("H2").Value = "2016-03-18 10:36:48"
